$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "reviews_count" column (column E). All columns to the right
# (reviews_average, latitude, longitude, is_permanently_closed, gmaps_link,
# latest_review_date) shift one position to the left.
$ws.Range("E1").EntireColumn.Delete()
